$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Rows.Item(6).Insert()
Write-Output $ws1.UsedRange.Address
